$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Mark Milestone I complete for the two rows that were previously blank.
$ws.Range("E6").Value = "I"
$ws.Range("F6").Value = "X"

$ws.Range("E10").Value = "I"
$ws.Range("F10").Value = "X"

# Mark Milestone II completion (X) for the carry-over rows.
$ws.Range("D91").Value = "X"
$ws.Range("D92").Value = "X"

$wb.Application.Calculate()

# Update the view: scroll so column C is leftmost and select D92.
$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 3
$ws.Range("D92").Select()
